$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "64.100.00"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "2.760.95"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "579.02"
$ws.Range("E5").Value = "  +0.31%  "
Set-TextValue $ws.Range("D6") "158.79"
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("E7").Value = "  +0.30%  "
Set-TextValue $ws.Range("D8") "0.608"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  -1.37%  "
Set-TextValue $ws.Range("D10") "5.72"
$ws.Range("E10").Value = "  -14.49%  "
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").Value = "3.249.10"
$ws.Range("E13").Value = "  +1.20%  "
Set-TextValue $ws.Range("D14") "26.99"
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("D15").Value = "63.767.43"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "2.761.96"
$ws.Range("E17").Value = "  +0.69%  "
Set-TextValue $ws.Range("D18") "12.16"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("E19").Value = "  +0.66%  "
Set-TextValue $ws.Range("D20") "361.41"
$ws.Range("E20").Value = "  +0.31%  "
Set-TextValue $ws.Range("D21") "6.85"
$ws.Range("E21").Value = "  -1.34%  "
Set-TextValue $ws.Range("D22") "0.551"
$ws.Range("E22").Value = "  +3.04%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("E24").Value = "  +0.16%  "
Set-TextValue $ws.Range("D25") "0.172"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "0.0₃0933"
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("E29").Value = "  -1.95%  "
Set-TextValue $ws.Range("D30") "7.05"
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("E31").Value = "  +0.83%  "
Set-TextValue $ws.Range("D32") "167.31"
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("E33").Value = "  -0.68%  "
Set-TextValue $ws.Range("D34") "4.96"
$ws.Range("E34").Value = "  +3.57%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("E36").Value = "  +2.40%  "
Set-TextValue $ws.Range("D37") "1.81"
$ws.Range("E37").Value = "  -0.22%  "
Set-TextValue $ws.Range("D38") "0.997"
$ws.Range("E38").Value = "  +0.04%  "
Set-TextValue $ws.Range("D39") "6.26"
$ws.Range("E39").Value = "  +12.24%  "
Set-TextValue $ws.Range("D40") "4.18"
$ws.Range("E40").Value = "  -0.91%  "
Set-TextValue $ws.Range("D41") "331.78"
$ws.Range("E41").Value = "  -3.70%  "
Set-TextValue $ws.Range("D42") "39.39"
$ws.Range("E42").Value = "  +0.24%  "
Set-TextValue $ws.Range("D43") "21.68"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("E45").Value = "  +0.85%  "
Set-TextValue $ws.Range("D46") "0.0259"
$ws.Range("E46").Value = "  +1.29%  "
Set-TextValue $ws.Range("D47") "0.636"
$ws.Range("E47").Value = "  -1.21%  "
Set-TextValue $ws.Range("D48") "136.11"
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("E51").Value = "  +0.70%  "
